# Updates cryptocurrency price (D) and 1h volume change (E) values
# to reflect the latest scrape, per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.709.38"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.851.02"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  -1.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.27"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4301"
$ws.Range("E7").Value = "  -2.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3745"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07340"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8773"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.59"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "1.840.10"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.736"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.434"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07136"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.00"
$ws.Range("E16").Value = "  +4.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008987"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.44"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").Value = "27.710.20"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.206"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").Value = "2.075.02"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.987"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.24"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.63"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.188"
$ws.Range("E28").Value = "  +9.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.369"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.93"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08937"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.229"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7782"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.550"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.931"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.012"
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.130"
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01983"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05342"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.302"
$ws.Range("E40").Value = "  +5.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.890"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1692"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5131"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.821"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.73"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.93"
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4779"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06468"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.012"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.689"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.841"
$ws.Range("E51").Value = "  -4.03%  "
